$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.353.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.57%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.478"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.71%  "
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.792.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.572.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.507"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.334.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0711"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.87%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("E26").Value = "  -2.89%  "
$ws.Range("E27").Value = "  -4.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.21%  "
$ws.Range("E30").Value = "  -6.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0462"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  -6.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.087.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  -4.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0151"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.781"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.494"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.34%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "93.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.89%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.755"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.706.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0108"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.25%  "
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.89%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.408"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
